$d = $word.ActiveDocument

# Remove "na uitleg A3 map" text (leaving the two preceding line breaks intact)
$d.Content.Find.Execute("na uitleg A3 map", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)

# Remove the entire "daarna uitleg Code structure/ flowcharts" text so the paragraph becomes empty
$d.Content.Find.Execute("daarna uitleg Code structure/ flowcharts", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
